# Update gh-pages to output generated at 1f05065
#
# For every sheet, column H ("是否有舞台（字符串匹配）") is removed. The
# columns to its right (I "Link", J "Cover") shift left to become H and I
# respectively, and the sheet's used range shrinks from J to I.
#
# In the two data sheets (展览 / sheet 1 and 全部类型 / sheet 4) the
# "想去人数" (F) and "最低票价" (G) values for rows 2-7 are also refreshed
# with newer scraped numbers, and G switches from a text value to a real
# number.

$wb = $excel.ActiveWorkbook

# Refreshed data for rows 2-7: 想去人数 (F) and 最低票价 (G)
$newData = @(
    @{ Row = 2; F = 2148; G = 54 },
    @{ Row = 3; F = 623;  G = 50 },
    @{ Row = 4; F = 1531; G = 60 },
    @{ Row = 5; F = 7243; G = 58 },
    @{ Row = 6; F = 177;  G = 50 },
    @{ Row = 7; F = 162;  G = 60 }
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Update the refreshed 想去人数 / 最低票价 values before the column
    # shift, while F/G still hold their original meaning.
    foreach ($entry in $newData) {
        $fCell = $ws.Cells.Item($entry.Row, 6)
        if (-not [string]::IsNullOrEmpty($fCell.Text)) {
            $ws.Cells.Item($entry.Row, 6).Value = $entry.F
            $ws.Cells.Item($entry.Row, 7).Value = $entry.G
        }
    }

    # Delete column H ("是否有舞台（字符串匹配）"); Link/Cover shift left.
    $ws.Columns.Item(8).Delete()
}
